# Weekly data refresh: a new record is inserted at row 200 (most recent
# week), and every existing record from row 200 downward shifts down by
# one row (the previously-oldest record, formerly row 315, ends up at
# the new row 316).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 200:315 down to 201:316, leaving a fresh blank row 200.
$ws.Range("A200").EntireRow.Insert()

# Populate the new row 200 with this week's record (same series as the
# row immediately below it, just a newer date/volume).
$ws.Range("A200").Value = 10
$ws.Range("B200").Value = "Vega Modelo de Temuco"
$ws.Range("C200").Value = "La Araucanía"
$ws.Range("D200").Value = 44603
$ws.Range("E200").Value = 9
$ws.Range("F200").Value = 100112037
$ws.Range("G200").Value = "Cebollín"
$ws.Range("H200").Value = "Sin especificar"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 95
$ws.Range("K200").Value = 7000
$ws.Range("L200").Value = 7000
$ws.Range("M200").Value = 7000
$ws.Range("N200").Value = "$/docena de paquetes"
$ws.Range("O200").Value = "Provincia de Cautín"
$ws.Range("P200").Value = 583
$ws.Range("Q200").Value = 12
$ws.Range("R200").Value = "Hortaliza"
